# Update SourceSinkModel results sheet values to reflect the new
# "shadow model" scaling (D1/L1 = D2/L2 --> D2 = ...)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = 3010.810655632908
$ws.Range("D3").Value  = 3010.810655632904
$ws.Range("D4").Value  = 3010.810655632904

$ws.Range("D6").Value  = 8602.31615895128
$ws.Range("D7").Value  = 8602.31615895128

$ws.Range("D9").Value  = 6876.196583938373
$ws.Range("D10").Value = 6876.196583938367
$ws.Range("D11").Value = 100
$ws.Range("D12").Value = 5676.196583938367
$ws.Range("D13").Value = 100

$ws.Range("D17").Value = 80000

$ws.Range("D19").Value = 87304.42145456493
$ws.Range("D20").Value = 87304.42145456493
$ws.Range("D21").Value = 1200

$ws.Range("D24").Value = 146534.5821035772
$ws.Range("D25").Value = 146534.582103577

$ws.Range("D28").Value = 2930.691642071566
$ws.Range("D29").Value = 2930.691642071566
$ws.Range("D30").Value = 146534.582103577

$ws.Range("D35").Value = 20217.60000000038
$ws.Range("D36").Value = 20217.60000000038

$ws.Range("D38").Value = -3776.765484902219
$ws.Range("D39").Value = -3776.765484902216

$ws.Range("D41").Value = 3776.765484902216
$ws.Range("D42").Value = 75535.30969804428
$ws.Range("D43").Value = 75535.30969804428
